# "Popup handling at profile management"
#
# The SignUp.xlsx fixture feeds the CreateAccount / Login Selenium
# scenarios. Each time the popup/profile-management flow is exercised the
# harness rotates the throw-away mailinator addresses used for the five
# "create account" rows plus the password used by the row-7 login check,
# so the same dummy users aren't reused across test runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateAccount")

# Rows 2-6: fresh Selenium mailinator addresses for each account type
# (Manufacturer, Retailer, Carrier, Designer, Individual Consumer).
$ws.Range("E2").Value = "SeleniumYcRj@mailinator.com"
$ws.Range("E3").Value = "SeleniumrFfh@mailinator.com"
$ws.Range("E4").Value = "SeleniummXTW@mailinator.com"
$ws.Range("E5").Value = "SeleniumbNGJ@mailinator.com"
$ws.Range("E6").Value = "SeleniumDjon@mailinator.com"

# Row 7: login-flow password rotated as well.
$ws.Range("F7").Value = "Automation4435!"

# Column E (email) widened slightly so the new values keep fitting without
# truncation, mirroring the auto-fit Excel applies after the edit.
$ws.Columns("E:E").ColumnWidth = 30
